$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 192; this shifts the existing rows 192-246
# down to 193-247 and extends the used range to A1:R247, matching the
# dimension change seen in the diff.
$ws.Rows("192:192").Insert()

# Populate the newly inserted row 192 with the new data point.
$ws.Range("A192").Value = 7
$ws.Range("B192").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C192").Value = "Ñuble"
$ws.Range("D192").Value = 44841
$ws.Range("E192").Value = 16
$ws.Range("F192").Value = 100112017
$ws.Range("G192").Value = "Apio"
$ws.Range("H192").Value = "Americana (o)"
$ws.Range("I192").Value = "Primera"
$ws.Range("J192").Value = 120
$ws.Range("K192").Value = 8500
$ws.Range("L192").Value = 9000
$ws.Range("M192").Value = 8750
$ws.Range("N192").Value = "$/docena de matas"
$ws.Range("O192").Value = "Provincia del Elquí"
$ws.Range("P192").Value = 1458
$ws.Range("Q192").Value = 6
$ws.Range("R192").Value = "Hortaliza"
